$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match data (xG_home, xG_away, goals_home, goals_away) for rows 10-15,
# which previously only had the home/away team names filled in (columns A-C).
# Values are written as text (matching the workbook's existing convention of
# storing these numeric-looking figures as shared strings rather than numbers).
$newData = @{
    10 = @{ D = "0.439998"; E = "1.64625";  F = "1"; G = "1" }
    11 = @{ D = "2.52836";  E = "1.34238";  F = "2"; G = "1" }
    12 = @{ D = "0.484378"; E = "3.60091";  F = "1"; G = "3" }
    13 = @{ D = "2.71039";  E = "0.902039"; F = "1"; G = "1" }
    14 = @{ D = "0.926619"; E = "2.80045";  F = "0"; G = "4" }
    15 = @{ D = "0.388088"; E = "1.53117";  F = "0"; G = "3" }
}

foreach ($row in $newData.Keys) {
    foreach ($col in "D", "E", "F", "G") {
        $text = $newData[$row][$col]
        $cell = $ws.Range("$col$row")
        # Prefix with an apostrophe so the numeric-looking text is stored as a
        # literal string (shared string) instead of being parsed into a number.
        $cell.Value = "'" + $text
        # Restore the default "Normal" style so no new/extra cell style is
        # introduced by the quote-prefix formatting.
        $cell.Style = "Normal"
    }
}
